# Fruta / hortaliza, semanal
#
# A new daily price-record row is inserted as row 31 (pushing the existing
# rows 31-111 down to 32-112), adding one more "Femacal de La Calera"
# Arándano (blue) observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 31; everything below shifts down one row.
$ws.Rows(31).Insert()

# Populate the newly inserted row 31 with the new observation.
$ws.Cells.Item(31, 1).Value  = 3
$ws.Cells.Item(31, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(31, 3).Value  = "Coquimbo"
$ws.Cells.Item(31, 4).Value  = 44481
$ws.Cells.Item(31, 5).Value  = 5
$ws.Cells.Item(31, 6).Value  = "Fruta"
$ws.Cells.Item(31, 7).Value  = 100101
$ws.Cells.Item(31, 8).Value  = "Berries"
$ws.Cells.Item(31, 9).Value  = 100101001
$ws.Cells.Item(31, 10).Value = "Arándano (blue)"
$ws.Cells.Item(31, 11).Value = "Sin especificar"
$ws.Cells.Item(31, 12).Value = "Primera"
$ws.Cells.Item(31, 13).Value = 95
$ws.Cells.Item(31, 14).Value = 11000
$ws.Cells.Item(31, 15).Value = 12000
$ws.Cells.Item(31, 16).Value = 11474
$ws.Cells.Item(31, 17).Value = "`$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(31, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(31, 19).Value = 7649
$ws.Cells.Item(31, 20).Value = 1.5

# Match the date formatting used by the rest of column D (YYYY-MM-DD HH:MM:SS).
$ws.Cells.Item(31, 4).NumberFormat = $ws.Cells.Item(32, 4).NumberFormat
